# Add a header row ("titles"/"ref"/"date") above the existing publication
# list on the "Sheet2" worksheet (the sorted/summarised table), shifting all
# data down by one row, then re-apply the existing sort (by column B,
# descending) over the new data range so sortState/sortCondition track the
# shifted range, and finally restore the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert a new blank row at the top; this pushes rows 1-255 down to 2-256
# and keeps every downstream style/value untouched.
$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Range("A1").Value = "titles"
$ws.Range("B1").Value = "ref"
$ws.Range("C1").Value = "date"

# Re-apply the descending sort on column B across the data (excluding the
# new header row) so the sheet's stored sort state matches the new extent.
$dataRange = $ws.Range("A2:C256")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B256"), 0, 2) | Out-Null
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Match the saved selection/active cell.
$ws.Range("E8").Select()
